$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") holds comma-separated lists of recorder names /
# emails. In a handful of rows the first two entries need to swap order
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
$map = @{
    "System, system, backup@backdoor.com" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
    "dnasr281@gmail.com, admin@admin.com" = "admin@admin.com, dnasr281@gmail.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey([string]$val)) {
        $cell.Value = $map[[string]$val]
    }
}
